$wb = $excel.ActiveWorkbook

# Insert a new worksheet as the first sheet (before the current first sheet),
# matching Excel's "right-click tab > Insert" behaviour: the new sheet becomes
# active/selected and the insertion point sheet loses focus.
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "watercolours"

# Populate the new sheet's data (3 columns x up to 6 rows).
$ws.Range("A1").Value = "Pink"
$ws.Range("B1").Value = "Purple"
$ws.Range("C1").Value = "Blue"

$ws.Range("A2").Value = "    Pastel Pink^"
$ws.Range("B2").Value = "Royal Purple"
$ws.Range("C2").Value = "Sea Blue"

$ws.Range("A3").Value = "  Rose Pink"
$ws.Range("B3").Value = "      Eggplant Purple"
$ws.Range("C3").Value = "Azure Blue      "

$ws.Range("A4").Value = "Very   Pink   "
$ws.Range("B4").Value = "Velvet Purple%`$"
$ws.Range("C4").Value = "Egg Blue       "

$ws.Range("A5").Value = "Rouge Pink"
$ws.Range("B5").Value = "Royal Purple"

$ws.Range("A6").Value = "Baby Pink"

# Size the columns to fit their (longest) contents, mirroring the author's
# manual column-width adjustment (AutoFit) after typing the data in.
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 16.833333333333332
$ws.Columns.Item(3).ColumnWidth = 12.5

# Put the cursor/selection where the author last left it on the new sheet.
$ws.Range("B20").Select() | Out-Null
